$d = $word.ActiveDocument

# The "Defendant ref:" merge field was rendering literally as
# "<<respondentExternalReference>>" because the party-index digit was
# missing from the field name (compare the sibling field just after it,
# "<<cs_{!isBlank(respondent2ExternalReference)}>>", and the claimant's
# field "<<applicant1ExternalReference>>" used elsewhere in the template).
# The fix is to insert the missing "1" so the field reads
# "<<respondent1ExternalReference>>".
#
# "respondentExternalReference>>" is unique in the document text stream
# (the neighbouring field already reads "respondent2ExternalReference>>"
# and the claimant field reads "applicantExternalReference>>", neither of
# which match this search string), so Find safely locates only the
# intended occurrence.
$found = $d.Content
$found.Find.Execute("respondentExternalReference>>", $true, $false, $false, $false, $false, `
                     $true, 1, $false, "", 0)

if ($found.Find.Found) {
    # Place the (collapsed) insertion point right after "respondent" - i.e.
    # between the "...nt" and "External..." text - and type the missing "1",
    # exactly as if a user had clicked there and pressed the "1" key.
    $insertionPoint = $d.Range($found.Start + 10, $found.Start + 10)
    $insertionPoint.InsertBefore("1")
}
